$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with former row 8 data (Mahle---Knecht)
$ws.Range("A2").Value = "Mahle---Knecht"
$ws.Range("B2").Value = "02943N0"
$ws.Range("C2").Value = "MZYJ"
$ws.Range("D2").Value = 52
$ws.Range("E2").Value = 1510

# Update row 3 with former row 5 data (Peugeot---Citroen / 82026)
$ws.Range("A3").Value = "Peugeot---Citroen"
$ws.Range("B3").Value = 82026
$ws.Range("C3").Value = "JFWU"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 309

# Remove the now-obsolete rows 4 through 8
$ws.Range("A4:E8").Delete()
